$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 13 de Octubre de 2020 a las 01:50"

# Update country data rows (refreshed COVID-19 statistics).
# Some rows also get a different country name because the countries
# swapped ranking position in the source data (sorted by Casos totales).

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 8036894
$ws.Range("C4").Value = 44896
$ws.Range("D4").Value = 5174828
$ws.Range("E4").Value = 2642056
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 315
$ws.Range("H4").Value = 220010

# Row 6: Brasil
$ws.Range("A6").Value = "Brasil"
$ws.Range("B6").Value = 5103408
$ws.Range("C6").Value = 8429
$ws.Range("D6").Value = 4495269
$ws.Range("E6").Value = 457430
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 203
$ws.Range("H6").Value = 150709

# Row 10: Argentina
$ws.Range("A10").Value = "Argentina"
$ws.Range("B10").Value = 903730
$ws.Range("C10").Value = 9524
$ws.Range("D10").Value = 732582
$ws.Range("E10").Value = 146962
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 318
$ws.Range("H10").Value = 24186

# Row 29: Canada
$ws.Range("A29").Value = "Canada"
$ws.Range("B29").Value = 182839
$ws.Range("C29").Value = 975
$ws.Range("D29").Value = 154258
$ws.Range("E29").Value = 18954
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 14
$ws.Range("H29").Value = 9627

# Row 38: Chequia
$ws.Range("A38").Value = "Chequia"
$ws.Range("B38").Value = 121421
$ws.Range("C38").Value = 4311
$ws.Range("D38").Value = 58200
$ws.Range("E38").Value = 62170
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 64
$ws.Range("H38").Value = 1051

# Row 39: Panama
$ws.Range("A39").Value = "Panama"
$ws.Range("B39").Value = 120802
$ws.Range("C39").Value = 489
$ws.Range("D39").Value = 96675
$ws.Range("E39").Value = 21625
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 11
$ws.Range("H39").Value = 2502

# Row 61: Nigeria
$ws.Range("A61").Value = "Nigeria"
$ws.Range("B61").Value = 60430
$ws.Range("C61").Value = 164
$ws.Range("D61").Value = 51943
$ws.Range("E61").Value = 7372
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 1115

# Row 67: Paraguay
$ws.Range("A67").Value = "Paraguay"
$ws.Range("B67").Value = 50344
$ws.Range("C67").Value = 669
$ws.Range("D67").Value = 32751
$ws.Range("E67").Value = 16497
$ws.Range("F67").Value = 0
$ws.Range("G67").Value = 19
$ws.Range("H67").Value = 1096

# Row 83: Australia
$ws.Range("A83").Value = "Australia"
$ws.Range("B83").Value = 27287
$ws.Range("C83").Value = 22
$ws.Range("D83").Value = 25013
$ws.Range("E83").Value = 1376
$ws.Range("F83").Value = 0
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 898

# Row 95: Noruega
$ws.Range("A95").Value = "Noruega"
$ws.Range("B95").Value = 15639
$ws.Range("C95").Value = 115
$ws.Range("D95").Value = 11863
$ws.Range("E95").Value = 3500
$ws.Range("F95").Value = 0
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = 276

# Row 117: Mauritania
$ws.Range("A117").Value = "Mauritania"
$ws.Range("B117").Value = 7554
$ws.Range("C117").Value = 4
$ws.Range("D117").Value = 7297
$ws.Range("E117").Value = 94
$ws.Range("F117").Value = 0
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 163

# Row 132: Surinam
$ws.Range("A132").Value = "Surinam"
$ws.Range("B132").Value = 5058
$ws.Range("C132").Value = 7
$ws.Range("D132").Value = 4862
$ws.Range("E132").Value = 89
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 107

# Row 149: Polinesia Francesa
$ws.Range("A149").Value = "Polinesia Francesa"
$ws.Range("B149").Value = 3251
$ws.Range("C149").Value = 497
$ws.Range("D149").Value = 2138
$ws.Range("E149").Value = 1102
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 1
$ws.Range("H149").Value = 11

# Row 150: Principado de Andorra
$ws.Range("A150").Value = "Principado de Andorra"
$ws.Range("B150").Value = 2995
$ws.Range("C150").Value = 299
$ws.Range("D150").Value = 1928
$ws.Range("E150").Value = 1010
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 40

# Row 151: Sudan del Sur
$ws.Range("A151").Value = "Sudan del Sur"
$ws.Range("B151").Value = 2787
$ws.Range("C151").Value = 10
$ws.Range("D151").Value = 1290
$ws.Range("E151").Value = 1442
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 55

# Row 152: Letonia
$ws.Range("A152").Value = "Letonia"
$ws.Range("B152").Value = 2765
$ws.Range("C152").Value = 95
$ws.Range("D152").Value = 1325
$ws.Range("E152").Value = 1400
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 40

# Row 156: Uruguay
$ws.Range("A156").Value = "Uruguay"
$ws.Range("B156").Value = 2313
$ws.Range("C156").Value = 19
$ws.Range("D156").Value = 1950
$ws.Range("E156").Value = 312
$ws.Range("F156").Value = 0
$ws.Range("G156").Value = 1
$ws.Range("H156").Value = 51

# Row 157: Sierra Leona
$ws.Range("A157").Value = "Sierra Leona"
$ws.Range("B157").Value = 2306
$ws.Range("C157").Value = 0
$ws.Range("D157").Value = 1736
$ws.Range("E157").Value = 498
$ws.Range("F157").Value = 0
$ws.Range("G157").Value = 0
$ws.Range("H157").Value = 72

# Row 158: Burkina Faso
$ws.Range("A158").Value = "Burkina Faso"
$ws.Range("B158").Value = 2294
$ws.Range("C158").Value = 23
$ws.Range("D158").Value = 1571
$ws.Range("E158").Value = 660
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 2
$ws.Range("H158").Value = 63
